# Update "想去人数" (F column) counts across the workbook's sheets.
# Sheets: 1=展览, 2=演出, 3=本地生活, 4=全部类型
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2887
$ws1.Range("F3").Value = 21155
$ws1.Range("F4").Value = 101
$ws1.Range("F5").Value = 2912
$ws1.Range("F6").Value = 808
$ws1.Range("F8").Value = 515
$ws1.Range("F9").Value = 765
$ws1.Range("F10").Value = 279
$ws1.Range("F12").Value = 72
$ws1.Range("F13").Value = 117
$ws1.Range("F14").Value = 518
$ws1.Range("F15").Value = 182
$ws1.Range("F16").Value = 268
$ws1.Range("F17").Value = 16
$ws1.Range("F18").Value = 422
$ws1.Range("F19").Value = 61
$ws1.Range("F22").Value = 34

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 32
$ws2.Range("F5").Value = 338
$ws2.Range("F6").Value = 142
$ws2.Range("F10").Value = 16
$ws2.Range("F14").Value = 153
$ws2.Range("F19").Value = 24
$ws2.Range("F22").Value = 40

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6131
$ws3.Range("F3").Value = 700
$ws3.Range("F4").Value = 695
$ws3.Range("F5").Value = 1608
$ws3.Range("F6").Value = 55

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6131
$ws4.Range("F3").Value = 700
$ws4.Range("F4").Value = 695
$ws4.Range("F5").Value = 1608
$ws4.Range("F6").Value = 2887
$ws4.Range("F7").Value = 21155
$ws4.Range("F9").Value = 32
$ws4.Range("F10").Value = 101
$ws4.Range("F12").Value = 338
$ws4.Range("F13").Value = 2912
$ws4.Range("F14").Value = 808
$ws4.Range("F15").Value = 142
$ws4.Range("F16").Value = 55
$ws4.Range("F18").Value = 516
$ws4.Range("F19").Value = 765
$ws4.Range("F20").Value = 279
$ws4.Range("F23").Value = 72
$ws4.Range("F26").Value = 117
$ws4.Range("F27").Value = 16
$ws4.Range("F29").Value = 518
$ws4.Range("F31").Value = 182
$ws4.Range("F33").Value = 268
$ws4.Range("F34").Value = 153
$ws4.Range("F35").Value = 153
$ws4.Range("F36").Value = 16
$ws4.Range("F37").Value = 422
$ws4.Range("F39").Value = 61
$ws4.Range("F44").Value = 34
$ws4.Range("F46").Value = 24
$ws4.Range("F49").Value = 40
